$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord coordinate values to whole numbers
$ws.Range("Q2").Value = 511515
$ws.Range("R2").Value = 6858547

# Clear the Starttid (Z2) and Sluttid (AB2) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
